$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "43.195.08"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "2.357.63"
$ws.Range("E3").Value = "  +5.42%  "

$ws.Range("E4").Value = "  +0.20%  "

Set-TextValue "D5" "233.67"
$ws.Range("E5").Value = "  +1.50%  "

Set-TextValue "D6" "0.652"
$ws.Range("E6").Value = "  +2.25%  "

Set-TextValue "D7" "72.46"
$ws.Range("E7").Value = "  +14.84%  "

$ws.Range("E8").Value = "  +0.00%  "

Set-TextValue "D9" "0.499"
$ws.Range("E9").Value = "  +13.67%  "

Set-TextValue "D10" "0.0977"
$ws.Range("E10").Value = "  +3.00%  "

Set-TextValue "D11" "27.30"
$ws.Range("E11").Value = "  -0.29%  "

Set-TextValue "D14" "16.29"
$ws.Range("E14").Value = "  +6.01%  "

Set-TextValue "D15" "6.29"
$ws.Range("E15").Value = "  +4.39%  "

Set-TextValue "D16" "0.865"
$ws.Range("E16").Value = "  +5.18%  "

$ws.Range("D17").Value = "2.361.80"
$ws.Range("E17").Value = "  +5.96%  "

$ws.Range("D18").Value = "43.217.54"
$ws.Range("E18").Value = "  +0.22%  "

Set-TextValue "D19" "0.0000100"
$ws.Range("E19").Value = "  +4.20%  "

$ws.Range("E20").Value = "  +5.01%  "

Set-TextValue "D21" "74.65"
$ws.Range("E21").Value = "  +2.74%  "

Set-TextValue "D22" "251.12"
$ws.Range("E22").Value = "  +2.36%  "

Set-TextValue "D23" "3.86"
$ws.Range("E23").Value = "  +5.84%  "

$ws.Range("E24").Value = "  +0.18%  "

Set-TextValue "D25" "2.46"
$ws.Range("E25").Value = "  +2.03%  "

Set-TextValue "D28" "22.48"
$ws.Range("E28").Value = "  +4.90%  "

Set-TextValue "D29" "172.59"
$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  +10.60%  "

$ws.Range("E31").Value = "  +4.39%  "

$ws.Range("E32").Value = "  +3.41%  "

Set-TextValue "D33" "5.01"
$ws.Range("E33").Value = "  +2.67%  "

Set-TextValue "D34" "0.0691"
$ws.Range("E34").Value = "  +3.00%  "

Set-TextValue "D35" "5.05"
$ws.Range("E35").Value = "  +4.14%  "

Set-TextValue "D37" "2.44"
$ws.Range("E37").Value = "  +7.57%  "

Set-TextValue "D38" "6.53"
$ws.Range("E38").Value = "  +4.33%  "

Set-TextValue "D39" "0.0255"
$ws.Range("E39").Value = "  +2.19%  "

Set-TextValue "D40" "19.17"
$ws.Range("E40").Value = "  +13.67%  "

$ws.Range("E41").Value = "  +0.06%  "

Set-TextValue "D42" "8.91"
$ws.Range("E42").Value = "  +3.40%  "

Set-TextValue "D46" "0.0960"
$ws.Range("E46").Value = "  +2.22%  "

$ws.Range("E47").Value = "  +2.88%  "

$ws.Range("D48").Value = "1.442.01"
$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("D49").Value = "2.584.71"
$ws.Range("E49").Value = "  +5.79%  "

$ws.Range("E50").Value = "  +1.18%  "

Set-TextValue "D51" "0.000202"
$ws.Range("E51").Value = "  -3.75%  "

# Rows 12/13 swap (TRON <-> WrappedliquidstakedEther2.0) with new data
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.713.45"
$ws.Range("E12").Value = "  +5.67%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D13" "0.106"
$ws.Range("E13").Value = "  +2.50%  "

# Rows 26/27 swap (Cosmos <-> Toncoin) with new data
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D26" "2.26"
$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D27" "10.03"
$ws.Range("E27").Value = "  +3.50%  "

# Rows 43/44/45 rotation (FTXToken/Aave/ARBITRUM) with new data
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "99.29"
$ws.Range("E43").Value = "  +3.38%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D44" "1.16"
$ws.Range("E44").Value = "  +9.90%  "

$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D45" "4.48"
$ws.Range("E45").Value = "  +0.61%  "

